$p = $ppt.ActivePresentation
$np = $p.Slides.Item(1).NotesPage
$cs = $np.ColorScheme
$cs.Item(3).RGB = 255
